$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 902.1539
$ws.Range("I6").Value = 675.1111
$ws.Range("K6").Value = 2025.3333
$ws.Range("M6").Value = -1913.3333
$ws.Range("H8").Value = 359.9
$ws.Range("I8").Value = 65.5
$ws.Range("K8").Value = 196.5
$ws.Range("M8").Value = -57.5
$ws.Range("H38").Value = 2627.9
$ws.Range("J38").Value = 6485
$ws.Range("L38").Value = 19455
$ws.Range("N38").Value = -20199
$ws.Range("H39").Value = 145.95238
$ws.Range("I39").Value = 121.6875
$ws.Range("J39").Value = 223.6
$ws.Range("K39").Value = 365.0625
$ws.Range("L39").Value = 670.8
$ws.Range("M39").Value = -69.0625
$ws.Range("N39").Value = -1262.8
$ws.Range("H76").Value = 299.5
$ws.Range("J76").Value = 299
$ws.Range("L76").Value = 299
$ws.Range("N76").Value = -929
$ws.Range("H79").Value = 299.5
$ws.Range("J79").Value = 299
$ws.Range("L79").Value = 299
$ws.Range("N79").Value = -2483
$ws.Range("H80").Value = 799.4286
$ws.Range("I80").Value = 649
$ws.Range("K80").Value = 1947
$ws.Range("M80").Value = -949
$ws.Range("H83").Value = 799.4286
$ws.Range("I83").Value = 649
$ws.Range("K83").Value = 5841
$ws.Range("M83").Value = -849
$ws.Range("H94").Value = 1300
$ws.Range("I94").Value = 1300
$ws.Range("K94").Value = 1300
$ws.Range("M94").Value = -849
$ws.Range("H116").Value = 3183
$ws.Range("I116").Value = 3175.3333
$ws.Range("K116").Value = 3175.3333
$ws.Range("M116").Value = 266.6667000000002
$ws.Range("H129").Value = 1735.6666
$ws.Range("I129").Value = 1397
$ws.Range("K129").Value = 4191
$ws.Range("M129").Value = 809
$ws.Range("H137").Value = 2899.8
$ws.Range("I137").Value = 1998.5
$ws.Range("J137").Value = 3038.4614
$ws.Range("K137").Value = 5995.5
$ws.Range("L137").Value = 9115.3842
$ws.Range("M137").Value = -3445.5
$ws.Range("N137").Value = -14215.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 943.2727
$ws.Range("I2").Value = 955.3333
$ws.Range("K2").Value = 955.3333
$ws.Range("M2").Value = -842.3333
$ws.Range("H32").Value = 489.42856
$ws.Range("I32").Value = 489.42856
$ws.Range("K32").Value = 489.42856
$ws.Range("M32").Value = -202.42856
$ws.Range("H37").Value = 46666.332
$ws.Range("J37").Value = 99999
$ws.Range("L37").Value = 99999
$ws.Range("N37").Value = -100545
$ws.Range("H45").Value = 2412.1538
$ws.Range("I45").Value = 1538
$ws.Range("K45").Value = 1538
$ws.Range("M45").Value = -1161
$ws.Range("H50").Value = 1132.6666
$ws.Range("J50").Value = 1625
$ws.Range("L50").Value = 1625
$ws.Range("N50").Value = -3053
$ws.Range("H61").Value = 3998.3333
$ws.Range("I61").Value = 3998.3333
$ws.Range("K61").Value = 3998.3333
$ws.Range("M61").Value = -3786.3333
$ws.Range("H116").Value = 943.2727
$ws.Range("I116").Value = 955.3333
$ws.Range("K116").Value = 955.3333
$ws.Range("M116").Value = 1338.6667
$ws.Range("H136").Value = 3998.3333
$ws.Range("I136").Value = 3998.3333
$ws.Range("K136").Value = 11994.9999
$ws.Range("M136").Value = -9444.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 943.2727
$ws.Range("I3").Value = 955.3333
$ws.Range("K3").Value = 955.3333
$ws.Range("M3").Value = -841.3333
$ws.Range("H20").Value = 1609.75
$ws.Range("I20").Value = 2175
$ws.Range("J20").Value = 1044.5
$ws.Range("K20").Value = 2175
$ws.Range("L20").Value = 1044.5
$ws.Range("M20").Value = -1928
$ws.Range("N20").Value = -1538.5
$ws.Range("H54").Value = 6577.6665
$ws.Range("I54").Value = 5366.5
$ws.Range("K54").Value = 5366.5
$ws.Range("M54").Value = -4882.5
$ws.Range("H80").Value = 322.8
$ws.Range("J80").Value = 434.4
$ws.Range("L80").Value = 434.4
$ws.Range("N80").Value = -2430.4
$ws.Range("H83").Value = 322.8
$ws.Range("J83").Value = 434.4
$ws.Range("L83").Value = 2172
$ws.Range("N83").Value = -12156
$ws.Range("H99").Value = 2211.3333
$ws.Range("J99").Value = 2499.5
$ws.Range("L99").Value = 2499.5
$ws.Range("N99").Value = -5495.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 68052.664
$ws.Range("I16").Value = 78181.92
$ws.Range("J16").Value = 2212.5
$ws.Range("K16").Value = 78181.92
$ws.Range("L16").Value = 2212.5
$ws.Range("M16").Value = -77894.92
$ws.Range("N16").Value = -2786.5
$ws.Range("H31").Value = 9359.5
$ws.Range("I31").Value = 3145
$ws.Range("J31").Value = 9837.538
$ws.Range("K31").Value = 3145
$ws.Range("L31").Value = 9837.538
$ws.Range("M31").Value = -2850
$ws.Range("N31").Value = -10427.538
$ws.Range("H34").Value = 9359.5
$ws.Range("I34").Value = 3145
$ws.Range("J34").Value = 9837.538
$ws.Range("K34").Value = 3145
$ws.Range("L34").Value = 9837.538
$ws.Range("M34").Value = -2943
$ws.Range("N34").Value = -10241.538
$ws.Range("H62").Value = 4830.6924
$ws.Range("J62").Value = 4900
$ws.Range("L62").Value = 4900
$ws.Range("N62").Value = -6148
$ws.Range("H65").Value = 4830.6924
$ws.Range("J65").Value = 4900
$ws.Range("L65").Value = 24500
$ws.Range("N65").Value = -30740
$ws.Range("H104").Value = 40285
$ws.Range("J104").Value = 40285
$ws.Range("L104").Value = 40285
$ws.Range("N104").Value = -45527
$ws.Range("H113").Value = 68052.664
$ws.Range("I113").Value = 78181.92
$ws.Range("J113").Value = 2212.5
$ws.Range("K113").Value = 78181.92
$ws.Range("L113").Value = 2212.5
$ws.Range("M113").Value = -76011.92
$ws.Range("N113").Value = -6552.5
$ws.Range("H132").Value = 2224.375
$ws.Range("I132").Value = 2320.7144
$ws.Range("K132").Value = 6962.1432
$ws.Range("M132").Value = -4432.1432
$ws.Range("H134").Value = 1700
$ws.Range("I134").Value = 1700
$ws.Range("K134").Value = 5100
$ws.Range("M134").Value = -2565

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 327.77777
$ws.Range("I44").Value = 162.5
$ws.Range("K44").Value = 487.5
$ws.Range("M44").Value = -89.5
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 45262
$ws.Range("J68").Value = 45262
$ws.Range("L68").Value = 45262
$ws.Range("N68").Value = -46884
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H71").Value = 45262
$ws.Range("J71").Value = 45262
$ws.Range("L71").Value = 135786
$ws.Range("N71").Value = -143898
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H113").Value = 7555
$ws.Range("I113").Value = 6495
$ws.Range("K113").Value = 6495
$ws.Range("M113").Value = -4325
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4500
$ws.Range("I132").Value = 4250
$ws.Range("K132").Value = 12750
$ws.Range("M132").Value = -10220
$ws.Range("H136").Value = 4992.6665
$ws.Range("I136").Value = 4989.5
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 14968.5
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -12418.5
$ws.Range("N136").Value = -20097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 448
$ws.Range("I55").Value = 448
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -171
$ws.Range("N55").ClearContents()
